# Update occupation classification code:
#  - insert a new "% of total people" column before "average debt per person"
#  - recompute average-debt-per-person (now column F) for every row
#  - reorder the occupation rows by descending total debt (column C)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- snapshot the existing data rows (rows 2-9, columns A-E) before touching anything ---
$rows = @()
for ($r = 2; $r -le 9; $r++) {
    $occupation = $ws.Cells.Item($r, 2).Value()
    $totalDebt = $ws.Cells.Item($r, 3).Value()
    $numPeople = $ws.Cells.Item($r, 4).Value()

    $rows += [PSCustomObject]@{
        Occupation = $occupation
        TotalDebt  = $totalDebt
        NumPeople  = $numPeople
    }
}

# sort the rows by total debt (column C) descending
$sorted = $rows | Sort-Object -Property TotalDebt -Descending

$totalPeople = 0
foreach ($row in $rows) { $totalPeople += $row.NumPeople }

# --- wipe the sheet so the shared-string table gets rebuilt in the new write order ---
$ws.Cells.Clear()

# --- header row ---
$ws.Range("B1").Value = "occupation"
$ws.Range("C1").Value = "6p_total"
$ws.Range("D1").Value = "# of people"
$ws.Range("E1").Value = "% of total people"
$ws.Range("F1").Value = "average debt per person"

$headerRange = $ws.Range("B1:F1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# --- write the reordered, recomputed data back out ---
$r = 2
foreach ($row in $sorted) {
    $index = $r - 2
    $avgDebt = $row.TotalDebt / $row.NumPeople
    $pctPeople = ($row.NumPeople / $totalPeople) * 100

    $aCell = $ws.Cells.Item($r, 1)
    $aCell.Value = $index
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1

    $ws.Cells.Item($r, 2).Value = $row.Occupation
    $ws.Cells.Item($r, 3).Value = $row.TotalDebt
    $ws.Cells.Item($r, 4).Value = $row.NumPeople
    $ws.Cells.Item($r, 5).Value = $pctPeople
    $ws.Cells.Item($r, 6).Value = $avgDebt

    $r++
}
